$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38; this shifts existing rows 38..97 down to 39..98
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record's data
$ws.Cells.Item(38, 1).Value2 = 5
$ws.Cells.Item(38, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(38, 3).Value2 = "Maule"
$ws.Cells.Item(38, 4).Value2 = 44797
$ws.Cells.Item(38, 5).Value2 = 7
$ws.Cells.Item(38, 6).Value2 = 100112013
$ws.Cells.Item(38, 7).Value2 = "Alcachofa"
$ws.Cells.Item(38, 8).Value2 = "Madrigal"
$ws.Cells.Item(38, 9).Value2 = "Primera"
$ws.Cells.Item(38, 10).Value2 = 400
$ws.Cells.Item(38, 11).Value2 = 13000
$ws.Cells.Item(38, 12).Value2 = 13000
$ws.Cells.Item(38, 13).Value2 = 13000
$ws.Cells.Item(38, 14).Value2 = "`$/caja 40 unidades"
$ws.Cells.Item(38, 15).Value2 = "Provincia del Elqu$([char]0x00ed)"
$ws.Cells.Item(38, 16).Value2 = 325
$ws.Cells.Item(38, 17).Value2 = 40
$ws.Cells.Item(38, 18).Value2 = "Hortaliza"

# Apply the same date-number-format style used by the other date cells in column D
$ws.Cells.Item(38, 4).NumberFormat = $ws.Cells.Item(39, 4).NumberFormat
